# "C suite.xlsx" - update the "Test Cases" sheet:
#   - Runmode column (C) switches from "Y" to "N" for every data row
#   - Results column (D) gets filled in with the outcome of each run
#   - selection/active cell moves from the Results column to the Runmode column
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("C2").Value = "N"
$ws.Range("D2").Value = "SKIP"

$ws.Range("C3").Value = "N"
$ws.Range("D3").Value = "SKIP"

$ws.Range("C4").Value = "N"
$ws.Range("D4").Value = "FAIL"

$ws.Range("C5").Value = "N"
$ws.Range("D5").Value = "SKIP"

$ws.Activate() | Out-Null
$ws.Range("C2:C5").Select() | Out-Null
